$d = $word.ActiveDocument

# Locate the paragraph that currently reads "Data, Information, Knowledge
# Model levels."
$findRng = $d.Content
$found = $findRng.Find.Execute("Data, Information, Knowledge Model levels.")

if (-not $found) {
    throw "Could not find anchor paragraph text"
}

$anchorStart = $findRng.Start

# Resolve the Paragraph index from the document-level Paragraphs
# collection (not from a derived sub-range), and keep re-resolving
# objects by fresh index lookups after each mutation below, since
# cached Paragraph/Range handles do not track shifting positions.
$count = $d.Paragraphs.Count
$anchorIndex = -1
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if (($p.Range.Start -le $anchorStart) -and ($p.Range.End -gt $anchorStart)) {
        $anchorIndex = $i
        break
    }
}

if ($anchorIndex -eq -1) {
    throw "Could not resolve anchor paragraph index"
}

# The existing (empty) spacer paragraph that already sits right after
# the anchor paragraph.
$spacerIndex = $anchorIndex + 1
$spacerPara = $d.Paragraphs.Item($spacerIndex)

# Insert two brand new paragraphs right before that existing spacer
# paragraph: first a blank paragraph, then a paragraph carrying the new
# body text -- matching the document's existing text / blank / text
# paragraph rhythm.
$insertPoint = $d.Range($spacerPara.Range.Start, $spacerPara.Range.Start)
$insertPoint.InsertParagraphBefore()

$spacerPara2 = $d.Paragraphs.Item($spacerIndex + 1)
$insertPoint2 = $d.Range($spacerPara2.Range.Start, $spacerPara2.Range.Start)
$insertPoint2.InsertParagraphBefore()

$newTextPara = $d.Paragraphs.Item($spacerIndex + 1)
$textRng = $d.Range($newTextPara.Range.Start, $newTextPara.Range.Start)
$textRng.InsertAfter("Reactive Resources: ID, IDOccurrence, Statement, Graph, Step, Messages, etc. Content Type Addressing: Graph Statements by patterns / RCV schema instances / roles / Kinds / SPARQL. Statement occurrences by position / role (Kind). IDOccurrences by contexts (schema / Kinds). IDs by occurrences (role / context schemas). Content Types: graph/set, graph/activation, graph/reference, statement/activation, occurrence/subject, producer/form (COST state exchange), etc. Super type / sub type functional transforms (Resources Function<Consumes, Produces>): addressing / query / traversal / augmentation steps.")

Write-Host "Inserted paragraphs after anchor paragraph index" $anchorIndex
